$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("M2").Value = 0.3007906666666667
$ws.Range("N2").Value = 0.902372
$ws.Range("O2").Value = 0.03537029821880876
$ws.Range("P2").Value = 0.03537029821880876
$ws.Range("Q2").Value = 2.260380197913333
$ws.Range("R2").Value = 20.34342178122
$ws.Range("S2").Value = 0.0345804616889741
$ws.Range("T2").Value = 0.0345804616889741

$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("O3").Value = 0.899334434508434
$ws.Range("P3").Value = 0.899334434508434
$ws.Range("S3").Value = 0.8792518447457264
$ws.Range("T3").Value = 0.8792518447457264

$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("M4").Value = 0.0008990000000000001
$ws.Range("N4").Value = 0.002697
$ws.Range("O4").Value = 0.0001057143775473167
$ws.Range("P4").Value = 0.0001057143775473167
$ws.Range("Q4").Value = 0.006755800705
$ws.Range("R4").Value = 0.060802206345
$ws.Range("S4").Value = 0.0001033537223840757
$ws.Range("T4").Value = 0.0001033537223840757

$ws.Range("I5").Value = 0.977669497583861
$ws.Range("J5").Value = 0.977669497583861
$ws.Range("M5").Value = 0.5528646666666667
$ws.Range("N5").Value = 1.658594
$ws.Range("O5").Value = 0.06501195117304938
$ws.Range("P5").Value = 0.06501195117304936
$ws.Range("Q5").Value = 4.154664632743333
$ws.Range("R5").Value = 37.39198169469
$ws.Range("S5").Value = 0.06356020164030168
$ws.Range("T5").Value = 0.06356020164030167

$ws.Range("I6").Value = 0.977669497583861
$ws.Range("J6").Value = 0.977669497583861
$ws.Range("M6").Value = 0.001510333333333333
$ws.Range("N6").Value = 0.004531
$ws.Range("O6").Value = 0.0001776017221605087
$ws.Range("P6").Value = 0.0001776017221605087
$ws.Range("Q6").Value = 0.01134984538166667
$ws.Range("R6").Value = 0.102148608435
$ws.Range("S6").Value = 0.000173635786474693
$ws.Range("T6").Value = 0.000173635786474693

$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.171642
$ws.Range("H7").Value = 0.514926
$ws.Range("I7").Value = 0.02233050241613897
$ws.Range("J7").Value = 0.02233050241613898
$ws.Range("M7").Value = 0.3007906666666667
$ws.Range("N7").Value = 0.902372
$ws.Range("O7").Value = 0.03537029821880876
$ws.Range("P7").Value = 0.03537029821880876
$ws.Range("Q7").Value = 0.05162831160799999
$ws.Range("R7").Value = 0.464654804472
$ws.Range("S7").Value = 0.0007898365298346651
$ws.Range("T7").Value = 0.0007898365298346652

$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.171642
$ws.Range("H8").Value = 0.514926
$ws.Range("I8").Value = 0.02233050241613897
$ws.Range("J8").Value = 0.02233050241613898
$ws.Range("O8").Value = 0.899334434508434
$ws.Range("P8").Value = 0.899334434508434
$ws.Range("Q8").Value = 1.312714926444
$ws.Range("R8").Value = 11.814434337996
$ws.Range("S8").Value = 0.02008258976270756
$ws.Range("T8").Value = 0.02008258976270757

$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.171642
$ws.Range("H9").Value = 0.514926
$ws.Range("I9").Value = 0.02233050241613897
$ws.Range("J9").Value = 0.02233050241613898
$ws.Range("M9").Value = 0.0008990000000000001
$ws.Range("N9").Value = 0.002697
$ws.Range("O9").Value = 0.0001057143775473167
$ws.Range("P9").Value = 0.0001057143775473167
$ws.Range("Q9").Value = 0.000154306158
$ws.Range("R9").Value = 0.001388755422
$ws.Range("S9").Value = 0.000002360655163240983
$ws.Range("T9").Value = 0.000002360655163240983

$ws.Range("E10").Value = 2.0
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.171642
$ws.Range("H10").Value = 0.514926
$ws.Range("I10").Value = 0.02233050241613897
$ws.Range("J10").Value = 0.02233050241613898
$ws.Range("M10").Value = 0.5528646666666667
$ws.Range("N10").Value = 1.658594
$ws.Range("O10").Value = 0.06501195117304938
$ws.Range("P10").Value = 0.06501195117304936
$ws.Range("Q10").Value = 0.09489479711599999
$ws.Range("R10").Value = 0.8540531740439999
$ws.Range("S10").Value = 0.001451749532747688
$ws.Range("T10").Value = 0.001451749532747688

$ws.Range("E11").Value = 2.0
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.171642
$ws.Range("H11").Value = 0.514926
$ws.Range("I11").Value = 0.02233050241613897
$ws.Range("J11").Value = 0.02233050241613898
$ws.Range("M11").Value = 0.001510333333333333
$ws.Range("N11").Value = 0.004531
$ws.Range("O11").Value = 0.0001776017221605087
$ws.Range("P11").Value = 0.0001776017221605087
$ws.Range("Q11").Value = 0.000259236634
$ws.Range("R11").Value = 0.002333129706
$ws.Range("S11").Value = 0.000003965935685815681
$ws.Range("T11").Value = 0.000003965935685815683
